$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Five monster-card effect cells get a wording change: "①在房间区" -> "①与玩家敌对"
# and "②在手牌" -> "②受玩家控制" (hostile-to-player / controlled-by-player deviation).
$ws.Range("E3").Replace("①在房间区", "①与玩家敌对")
$ws.Range("E3").Replace("②在手牌", "②受玩家控制")

$ws.Range("E4").Replace("①在房间区", "①与玩家敌对")
$ws.Range("E4").Replace("②在手牌", "②受玩家控制")

$ws.Range("E5").Replace("①在房间区", "①与玩家敌对")
$ws.Range("E5").Replace("②在手牌", "②受玩家控制")

$ws.Range("E6").Replace("①在房间区", "①与玩家敌对")
$ws.Range("E6").Replace("②在手牌", "②受玩家控制")

$ws.Range("E14").Replace("①在房间区", "①与玩家敌对")
$ws.Range("E14").Replace("②在手牌", "②受玩家控制")

# End on E14 to match the author's final selection/scroll position.
$ws.Range("A14").Select()
$ws.Range("E14").Select()
